$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.993.29"
$ws.Range("E2").Value = "  +1.86%  "
$ws.Range("D3").Value = "1.907.13"
$ws.Range("E3").Value = "  +2.20%  "
$ws.Range("E4").Value = "  -0.89%  "
$ws.Range("D5").Value = "'315.49"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.20%  "
$ws.Range("E6").Value = "  -0.78%  "
$ws.Range("D7").Value = "'0.4816"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.85%  "
$ws.Range("D8").Value = "'0.3799"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.58%  "
$ws.Range("D9").Value = "'0.07361"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.33%  "
$ws.Range("D10").Value = "'0.9318"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.44%  "
$ws.Range("D11").Value = "'20.78"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.51%  "
$ws.Range("D12").Value = "'0.07754"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.09%  "
$ws.Range("D13").Value = "1.951.66"
$ws.Range("E13").Value = "  +3.54%  "
$ws.Range("D14").Value = "'5.498"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.06%  "
$ws.Range("D15").Value = "'6.631"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.23%  "
$ws.Range("D16").Value = "'91.71"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.39%  "
$ws.Range("D17").Value = "'1.005"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.86%  "
$ws.Range("D18").Value = "'0.000008823"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.83%  "
$ws.Range("D19").Value = "'1.003"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.76%  "
$ws.Range("D20").Value = "28.024.70"
$ws.Range("E20").Value = "  +1.69%  "
$ws.Range("D21").Value = "'14.78"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.38%  "
$ws.Range("D22").Value = "'5.164"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.85%  "
$ws.Range("D23").Value = "2.158.55"
$ws.Range("E23").Value = "  +1.31%  "
$ws.Range("D24").Value = "'10.87"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.65%  "
$ws.Range("D25").Value = "'155.50"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.69%  "
$ws.Range("E26").Value = "  -1.21%  "
$ws.Range("E27").Value = "  -0.01%  "
$ws.Range("D28").Value = "'2.130"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +5.42%  "
$ws.Range("D29").Value = "'116.98"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.26%  "
$ws.Range("D30").Value = "'4.955"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.51%  "
$ws.Range("D31").Value = "'0.08943"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.50%  "
$ws.Range("D32").Value = "'3.297"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.97%  "
$ws.Range("E33").Value = "  +3.28%  "
$ws.Range("D34").Value = "'0.7744"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.03%  "
$ws.Range("D35").Value = "'4.678"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.42%  "
$ws.Range("D36").Value = "'2.636"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.43%  "
$ws.Range("D37").Value = "'0.02056"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.11%  "
$ws.Range("D38").Value = "'1.110"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.96%  "
$ws.Range("D39").Value = "'0.05309"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.70%  "
$ws.Range("D40").Value = "'0.5485"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.12%  "
$ws.Range("D41").Value = "'2.994"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.04%  "
$ws.Range("D42").Value = "'7.019"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.01%  "
$ws.Range("D43").Value = "'0.1526"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.08%  "
$ws.Range("D44").Value = "'8.479"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.30%  "
$ws.Range("D45").Value = "'10.68"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.43%  "
$ws.Range("D46").Value = "'0.4828"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.42%  "
$ws.Range("D47").Value = "'108.08"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.88%  "
$ws.Range("E48").Value = "  -0.81%  "
$ws.Range("D49").Value = "'1.650"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.16%  "
$ws.Range("D50").Value = "'67.88"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.60%  "
$ws.Range("D51").Value = "'0.06070"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.18%  "
